# "Generate Report for Handoff" - refresh the handoff-status report:
#   - Overview sheet: refresh the "Latest HO Xliff Generate Date" for the
#     00a5926a-... file (rows 4-7), reflecting the new handoff timestamp.
#   - zh-cn / de-de sheets: promote rows 4-7 (the 00a5926a/10aa8b00/
#     7cc3fefa/ac12d43b files) from priority "low" to "ht", and refresh the
#     zh-cn "Latest Handoff Datetime" to the new handoff run's timestamp.

$wb = $excel.ActiveWorkbook
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

foreach ($r in 4..7) {
    # Priority low -> ht for both locales
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # New handoff generation timestamp, shared by Overview + de-de (same
    # underlying value, so it round-trips as a single shared string)
    $overview.Range("G$r").Value = "2016-08-26 08:35:33"
    $dede.Range("H$r").Value     = "2016-08-26 08:35:33"

    # zh-cn's own Latest Handoff Datetime moves forward too
    $zhcn.Range("H$r").Value = "2016-08-26 08:35:29"
}
